$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the bordered/centered numbering style from A22 down through A50
# (mirrors rows 23-50 being new entries in column A/C).
$ws.Range("A22").Copy()
$ws.Range("A23:A50").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate columns A-F for rows 2..50 with the refreshed screener data.
# Column A keeps its running 0-based index; D and F are always blank.

$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = "NSE:CYBERMEDIA"
$ws.Cells.Item(2, 3).Value = "NSE:ACE"
$ws.Cells.Item(2, 4).ClearContents()
$ws.Cells.Item(2, 5).Value = "NSE:ABCAPITAL"
$ws.Cells.Item(2, 6).ClearContents()

$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = "NSE:DOLATALGO"
$ws.Cells.Item(3, 3).Value = "NSE:ADFFOODS"
$ws.Cells.Item(3, 4).ClearContents()
$ws.Cells.Item(3, 5).Value = "NSE:BAJFINANCE"
$ws.Cells.Item(3, 6).ClearContents()

$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(4, 2).Value = "NSE:DONEAR"
$ws.Cells.Item(4, 3).Value = "NSE:AIAENG"
$ws.Cells.Item(4, 4).ClearContents()
$ws.Cells.Item(4, 5).Value = "NSE:BANDHANBNK"
$ws.Cells.Item(4, 6).ClearContents()

$ws.Cells.Item(5, 1).Value = 3
$ws.Cells.Item(5, 2).Value = "NSE:GILLETTE"
$ws.Cells.Item(5, 3).Value = "NSE:ALOKINDS"
$ws.Cells.Item(5, 4).ClearContents()
$ws.Cells.Item(5, 5).Value = "NSE:CANBK"
$ws.Cells.Item(5, 6).ClearContents()

$ws.Cells.Item(6, 1).Value = 4
$ws.Cells.Item(6, 2).Value = "NSE:GLAXO"
$ws.Cells.Item(6, 3).Value = "NSE:AMDIND"
$ws.Cells.Item(6, 4).ClearContents()
$ws.Cells.Item(6, 5).Value = "NSE:HAL"
$ws.Cells.Item(6, 6).ClearContents()

$ws.Cells.Item(7, 1).Value = 5
$ws.Cells.Item(7, 2).Value = "NSE:HARSHA"
$ws.Cells.Item(7, 3).Value = "NSE:ANDHRSUGAR"
$ws.Cells.Item(7, 4).ClearContents()
$ws.Cells.Item(7, 5).Value = "NSE:ICICIBANK"
$ws.Cells.Item(7, 6).ClearContents()

$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 2).Value = "NSE:HINDZINC"
$ws.Cells.Item(8, 3).Value = "NSE:ASAHIINDIA"
$ws.Cells.Item(8, 4).ClearContents()
$ws.Cells.Item(8, 5).Value = "NSE:INDHOTEL"
$ws.Cells.Item(8, 6).ClearContents()

$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = "NSE:ICEMAKE"
$ws.Cells.Item(9, 3).Value = "NSE:ASAHISONG"
$ws.Cells.Item(9, 4).ClearContents()
$ws.Cells.Item(9, 5).Value = "NSE:INDUSINDBK"
$ws.Cells.Item(9, 6).ClearContents()

$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 2).Value = "NSE:MASTEK"
$ws.Cells.Item(10, 3).Value = "NSE:BAJAJHIND"
$ws.Cells.Item(10, 4).ClearContents()
$ws.Cells.Item(10, 5).Value = "NSE:LICHSGFIN"
$ws.Cells.Item(10, 6).ClearContents()

$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(11, 2).Value = "NSE:MSPL"
$ws.Cells.Item(11, 3).Value = "NSE:BHEL"
$ws.Cells.Item(11, 4).ClearContents()
$ws.Cells.Item(11, 5).Value = "NSE:OFSS"
$ws.Cells.Item(11, 6).ClearContents()

$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(12, 2).Value = "NSE:MUTHOOTFIN"
$ws.Cells.Item(12, 3).Value = "NSE:CENTURYTEX"
$ws.Cells.Item(12, 4).ClearContents()
$ws.Cells.Item(12, 5).ClearContents()
$ws.Cells.Item(12, 6).ClearContents()

$ws.Cells.Item(13, 1).Value = 11
$ws.Cells.Item(13, 2).Value = "NSE:NAZARA"
$ws.Cells.Item(13, 3).Value = "NSE:COMPUSOFT"
$ws.Cells.Item(13, 4).ClearContents()
$ws.Cells.Item(13, 5).ClearContents()
$ws.Cells.Item(13, 6).ClearContents()

$ws.Cells.Item(14, 1).Value = 12
$ws.Cells.Item(14, 2).Value = "NSE:PRITIKAUTO"
$ws.Cells.Item(14, 3).Value = "NSE:DWARKESH"
$ws.Cells.Item(14, 4).ClearContents()
$ws.Cells.Item(14, 5).ClearContents()
$ws.Cells.Item(14, 6).ClearContents()

$ws.Cells.Item(15, 1).Value = 13
$ws.Cells.Item(15, 2).Value = "NSE:RADHIKAJWE"
$ws.Cells.Item(15, 3).Value = "NSE:FACT"
$ws.Cells.Item(15, 4).ClearContents()
$ws.Cells.Item(15, 5).ClearContents()
$ws.Cells.Item(15, 6).ClearContents()

$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = "NSE:ROSSARI"
$ws.Cells.Item(16, 3).Value = "NSE:FCL"
$ws.Cells.Item(16, 4).ClearContents()
$ws.Cells.Item(16, 5).ClearContents()
$ws.Cells.Item(16, 6).ClearContents()

$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 2).ClearContents()
$ws.Cells.Item(17, 3).Value = "NSE:FINEORG"
$ws.Cells.Item(17, 4).ClearContents()
$ws.Cells.Item(17, 5).ClearContents()
$ws.Cells.Item(17, 6).ClearContents()

$ws.Cells.Item(18, 1).Value = 16
$ws.Cells.Item(18, 2).ClearContents()
$ws.Cells.Item(18, 3).Value = "NSE:FINPIPE"
$ws.Cells.Item(18, 4).ClearContents()
$ws.Cells.Item(18, 5).ClearContents()
$ws.Cells.Item(18, 6).ClearContents()

$ws.Cells.Item(19, 1).Value = 17
$ws.Cells.Item(19, 2).ClearContents()
$ws.Cells.Item(19, 3).Value = "NSE:GATEWAY"
$ws.Cells.Item(19, 4).ClearContents()
$ws.Cells.Item(19, 5).ClearContents()
$ws.Cells.Item(19, 6).ClearContents()

$ws.Cells.Item(20, 1).Value = 18
$ws.Cells.Item(20, 2).ClearContents()
$ws.Cells.Item(20, 3).Value = "NSE:GREAVESCOT"
$ws.Cells.Item(20, 4).ClearContents()
$ws.Cells.Item(20, 5).ClearContents()
$ws.Cells.Item(20, 6).ClearContents()

$ws.Cells.Item(21, 1).Value = 19
$ws.Cells.Item(21, 2).ClearContents()
$ws.Cells.Item(21, 3).Value = "NSE:GRINDWELL"
$ws.Cells.Item(21, 4).ClearContents()
$ws.Cells.Item(21, 5).ClearContents()
$ws.Cells.Item(21, 6).ClearContents()

$ws.Cells.Item(22, 1).Value = 20
$ws.Cells.Item(22, 2).ClearContents()
$ws.Cells.Item(22, 3).Value = "NSE:GRSE"
$ws.Cells.Item(22, 4).ClearContents()
$ws.Cells.Item(22, 5).ClearContents()
$ws.Cells.Item(22, 6).ClearContents()

$ws.Cells.Item(23, 1).Value = 21
$ws.Cells.Item(23, 2).ClearContents()
$ws.Cells.Item(23, 3).Value = "NSE:HEIDELBERG"
$ws.Cells.Item(23, 4).ClearContents()
$ws.Cells.Item(23, 5).ClearContents()
$ws.Cells.Item(23, 6).ClearContents()

$ws.Cells.Item(24, 1).Value = 22
$ws.Cells.Item(24, 2).ClearContents()
$ws.Cells.Item(24, 3).Value = "NSE:HILTON"
$ws.Cells.Item(24, 4).ClearContents()
$ws.Cells.Item(24, 5).ClearContents()
$ws.Cells.Item(24, 6).ClearContents()

$ws.Cells.Item(25, 1).Value = 23
$ws.Cells.Item(25, 2).ClearContents()
$ws.Cells.Item(25, 3).Value = "NSE:IDEAFORGE"
$ws.Cells.Item(25, 4).ClearContents()
$ws.Cells.Item(25, 5).ClearContents()
$ws.Cells.Item(25, 6).ClearContents()

$ws.Cells.Item(26, 1).Value = 24
$ws.Cells.Item(26, 2).ClearContents()
$ws.Cells.Item(26, 3).Value = "NSE:INDIAMART"
$ws.Cells.Item(26, 4).ClearContents()
$ws.Cells.Item(26, 5).ClearContents()
$ws.Cells.Item(26, 6).ClearContents()

$ws.Cells.Item(27, 1).Value = 25
$ws.Cells.Item(27, 2).ClearContents()
$ws.Cells.Item(27, 3).Value = "NSE:INDUSTOWER"
$ws.Cells.Item(27, 4).ClearContents()
$ws.Cells.Item(27, 5).ClearContents()
$ws.Cells.Item(27, 6).ClearContents()

$ws.Cells.Item(28, 1).Value = 26
$ws.Cells.Item(28, 2).ClearContents()
$ws.Cells.Item(28, 3).Value = "NSE:ITI"
$ws.Cells.Item(28, 4).ClearContents()
$ws.Cells.Item(28, 5).ClearContents()
$ws.Cells.Item(28, 6).ClearContents()

$ws.Cells.Item(29, 1).Value = 27
$ws.Cells.Item(29, 2).ClearContents()
$ws.Cells.Item(29, 3).Value = "NSE:JCHAC"
$ws.Cells.Item(29, 4).ClearContents()
$ws.Cells.Item(29, 5).ClearContents()
$ws.Cells.Item(29, 6).ClearContents()

$ws.Cells.Item(30, 1).Value = 28
$ws.Cells.Item(30, 2).ClearContents()
$ws.Cells.Item(30, 3).Value = "NSE:JSL"
$ws.Cells.Item(30, 4).ClearContents()
$ws.Cells.Item(30, 5).ClearContents()
$ws.Cells.Item(30, 6).ClearContents()

$ws.Cells.Item(31, 1).Value = 29
$ws.Cells.Item(31, 2).ClearContents()
$ws.Cells.Item(31, 3).Value = "NSE:JYOTHYLAB"
$ws.Cells.Item(31, 4).ClearContents()
$ws.Cells.Item(31, 5).ClearContents()
$ws.Cells.Item(31, 6).ClearContents()

$ws.Cells.Item(32, 1).Value = 30
$ws.Cells.Item(32, 2).ClearContents()
$ws.Cells.Item(32, 3).Value = "NSE:KABRAEXTRU"
$ws.Cells.Item(32, 4).ClearContents()
$ws.Cells.Item(32, 5).ClearContents()
$ws.Cells.Item(32, 6).ClearContents()

$ws.Cells.Item(33, 1).Value = 31
$ws.Cells.Item(33, 2).ClearContents()
$ws.Cells.Item(33, 3).Value = "NSE:KAJARIACER"
$ws.Cells.Item(33, 4).ClearContents()
$ws.Cells.Item(33, 5).ClearContents()
$ws.Cells.Item(33, 6).ClearContents()

$ws.Cells.Item(34, 1).Value = 32
$ws.Cells.Item(34, 2).ClearContents()
$ws.Cells.Item(34, 3).Value = "NSE:KIMS"
$ws.Cells.Item(34, 4).ClearContents()
$ws.Cells.Item(34, 5).ClearContents()
$ws.Cells.Item(34, 6).ClearContents()

$ws.Cells.Item(35, 1).Value = 33
$ws.Cells.Item(35, 2).ClearContents()
$ws.Cells.Item(35, 3).Value = "NSE:KNRCON"
$ws.Cells.Item(35, 4).ClearContents()
$ws.Cells.Item(35, 5).ClearContents()
$ws.Cells.Item(35, 6).ClearContents()

$ws.Cells.Item(36, 1).Value = 34
$ws.Cells.Item(36, 2).ClearContents()
$ws.Cells.Item(36, 3).Value = "NSE:MAITHANALL"
$ws.Cells.Item(36, 4).ClearContents()
$ws.Cells.Item(36, 5).ClearContents()
$ws.Cells.Item(36, 6).ClearContents()

$ws.Cells.Item(37, 1).Value = 35
$ws.Cells.Item(37, 2).ClearContents()
$ws.Cells.Item(37, 3).Value = "NSE:MANGCHEFER"
$ws.Cells.Item(37, 4).ClearContents()
$ws.Cells.Item(37, 5).ClearContents()
$ws.Cells.Item(37, 6).ClearContents()

$ws.Cells.Item(38, 1).Value = 36
$ws.Cells.Item(38, 2).ClearContents()
$ws.Cells.Item(38, 3).Value = "NSE:MBAPL"
$ws.Cells.Item(38, 4).ClearContents()
$ws.Cells.Item(38, 5).ClearContents()
$ws.Cells.Item(38, 6).ClearContents()

$ws.Cells.Item(39, 1).Value = 37
$ws.Cells.Item(39, 2).ClearContents()
$ws.Cells.Item(39, 3).Value = "NSE:MMTC"
$ws.Cells.Item(39, 4).ClearContents()
$ws.Cells.Item(39, 5).ClearContents()
$ws.Cells.Item(39, 6).ClearContents()

$ws.Cells.Item(40, 1).Value = 38
$ws.Cells.Item(40, 2).ClearContents()
$ws.Cells.Item(40, 3).Value = "NSE:MOIL"
$ws.Cells.Item(40, 4).ClearContents()
$ws.Cells.Item(40, 5).ClearContents()
$ws.Cells.Item(40, 6).ClearContents()

$ws.Cells.Item(41, 1).Value = 39
$ws.Cells.Item(41, 2).ClearContents()
$ws.Cells.Item(41, 3).Value = "NSE:NAGAFERT"
$ws.Cells.Item(41, 4).ClearContents()
$ws.Cells.Item(41, 5).ClearContents()
$ws.Cells.Item(41, 6).ClearContents()

$ws.Cells.Item(42, 1).Value = 40
$ws.Cells.Item(42, 2).ClearContents()
$ws.Cells.Item(42, 3).Value = "NSE:NILKAMAL"
$ws.Cells.Item(42, 4).ClearContents()
$ws.Cells.Item(42, 5).ClearContents()
$ws.Cells.Item(42, 6).ClearContents()

$ws.Cells.Item(43, 1).Value = 41
$ws.Cells.Item(43, 2).ClearContents()
$ws.Cells.Item(43, 3).Value = "NSE:ONWARDTEC"
$ws.Cells.Item(43, 4).ClearContents()
$ws.Cells.Item(43, 5).ClearContents()
$ws.Cells.Item(43, 6).ClearContents()

$ws.Cells.Item(44, 1).Value = 42
$ws.Cells.Item(44, 2).ClearContents()
$ws.Cells.Item(44, 3).Value = "NSE:ORISSAMINE"
$ws.Cells.Item(44, 4).ClearContents()
$ws.Cells.Item(44, 5).ClearContents()
$ws.Cells.Item(44, 6).ClearContents()

$ws.Cells.Item(45, 1).Value = 43
$ws.Cells.Item(45, 2).ClearContents()
$ws.Cells.Item(45, 3).Value = "NSE:PRAJIND"
$ws.Cells.Item(45, 4).ClearContents()
$ws.Cells.Item(45, 5).ClearContents()
$ws.Cells.Item(45, 6).ClearContents()

$ws.Cells.Item(46, 1).Value = 44
$ws.Cells.Item(46, 2).ClearContents()
$ws.Cells.Item(46, 3).Value = "NSE:RATEGAIN"
$ws.Cells.Item(46, 4).ClearContents()
$ws.Cells.Item(46, 5).ClearContents()
$ws.Cells.Item(46, 6).ClearContents()

$ws.Cells.Item(47, 1).Value = 45
$ws.Cells.Item(47, 2).ClearContents()
$ws.Cells.Item(47, 3).Value = "NSE:RCF"
$ws.Cells.Item(47, 4).ClearContents()
$ws.Cells.Item(47, 5).ClearContents()
$ws.Cells.Item(47, 6).ClearContents()

$ws.Cells.Item(48, 1).Value = 46
$ws.Cells.Item(48, 2).ClearContents()
$ws.Cells.Item(48, 3).Value = "NSE:RTNINDIA"
$ws.Cells.Item(48, 4).ClearContents()
$ws.Cells.Item(48, 5).ClearContents()
$ws.Cells.Item(48, 6).ClearContents()

$ws.Cells.Item(49, 1).Value = 47
$ws.Cells.Item(49, 2).ClearContents()
$ws.Cells.Item(49, 3).Value = "NSE:SAGCEM"
$ws.Cells.Item(49, 4).ClearContents()
$ws.Cells.Item(49, 5).ClearContents()
$ws.Cells.Item(49, 6).ClearContents()

$ws.Cells.Item(50, 1).Value = 48
$ws.Cells.Item(50, 2).ClearContents()
$ws.Cells.Item(50, 3).Value = "NSE:SAKHTISUG"
$ws.Cells.Item(50, 4).ClearContents()
$ws.Cells.Item(50, 5).ClearContents()
$ws.Cells.Item(50, 6).ClearContents()
